$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.996.18"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "'2.496.43"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'536.70"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "'143.09"
$ws.Range("E6").Value = "  -3.61%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'2.529.77"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("D10").Value = "'0.0993"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "'5.46"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").Value = "'0.350"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "'2.973.40"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "'23.55"
$ws.Range("E15").Value = "  -3.49%  "
$ws.Range("D16").Value = "'58.980.96"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "'2.522.48"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'11.22"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "'4.27"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").Value = "'322.48"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("D23").Value = "'5.76"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'62.00"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").Value = "'0.436"
$ws.Range("E25").Value = "  -8.09%  "
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").Value = "'2.627.69"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").Value = "'7.73"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").Value = "'6.71"
$ws.Range("E30").Value = "  -4.21%  "
$ws.Range("D31").Value = "'0.0₃0770"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("D33").Value = "'1.18"
$ws.Range("E33").Value = "  -7.84%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'158.27"
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").Value = "'1.43"
$ws.Range("E36").Value = "  +5.70%  "
$ws.Range("D37").Value = "'18.56"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("D38").Value = "'4.34"
$ws.Range("E38").Value = "  -5.94%  "
$ws.Range("E39").Value = "  -5.72%  "
$ws.Range("D40").Value = "'5.62"
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("D41").Value = "'36.87"
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "'300.13"
$ws.Range("E42").Value = "  -4.02%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "'3.64"
$ws.Range("E43").Value = "  -3.70%  "
$ws.Range("B44").Value = "SuiNetwork"
$ws.Range("C44").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D44").Value = "'0.811"
$ws.Range("E44").Value = "  -6.14%  "
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("E46").Value = "  +3.31%  "
$ws.Range("D47").Value = "'10.77"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").Value = "'125.47"
$ws.Range("E48").Value = "  +4.84%  "
$ws.Range("D49").Value = "'0.0929"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "'18.66"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'0.0513"
$ws.Range("E51").Value = "  -2.15%  "
